$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-29 Friday" "2024-03-30 Saturday"

Replace-Text "19×16=" "70×69="
Replace-Text "58×20=" "94×91="
Replace-Text "74×70=" "48×50="
Replace-Text "14×75=" "63×79="
Replace-Text "37×65=" "34×45="
Replace-Text "31×14=" "38×98="
Replace-Text "44×41=" "17×65="
Replace-Text "13×64=" "56×63="
Replace-Text "44×70=" "63×55="
Replace-Text "43×56=" "72×46="
Replace-Text "74×59=" "28×72="
Replace-Text "19×19=" "92×31="
Replace-Text "37×61=" "68×51="
Replace-Text "91×40=" "86×88="
Replace-Text "17×51=" "32×94="
Replace-Text "33×54=" "43×66="
Replace-Text "94×67=" "50×63="
Replace-Text "62×31=" "31×87="
Replace-Text "23×66=" "69×59="
Replace-Text "16×53=" "74×70="
Replace-Text "25×58=" "39×71="
Replace-Text "29×66=" "82×72="
Replace-Text "92×64=" "34×77="
Replace-Text "69×25=" "29×63="
Replace-Text "38×64=" "73×32="
